$wb = $excel.ActiveWorkbook

# Rename sheets (tab names)
$wb.Worksheets.Item("GNG_TO-16504777959355574").Name = "GNG_TO-16509960845067885"
$wb.Worksheets.Item("NB_TO-16504777979295547").Name = "NB_TO-1650996086995259"
$wb.Worksheets.Item("RS_TO-16504777979305549").Name = "RS_TO-1650996086995259"
$wb.Worksheets.Item("TOL_TO-16504777979925542").Name = "TOL_TO-16509960870432591"
$wb.Worksheets.Item("vSAT_TO-16504777980555944").Name = "vSAT_TO-16509960871165123"

# Sheet1: GNG_TO
$ws1 = $wb.Worksheets.Item("GNG_TO-16509960845067885")
$ws1.Range("B2").Value = "go_stims-16509960844747891.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960844908254.csv"
$ws1.Range("B4").Value = "go_stims-16509960844908254.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960845067885.csv"

# Sheet2: NB_TO
$ws2 = $wb.Worksheets.Item("NB_TO-1650996086995259")
$ws2.Range("B2").Value = "OB-16509960854992585.csv"
$ws2.Range("B3").Value = "ZB-match_8-165099608462679.csv"
$ws2.Range("B4").Value = "ZB-match_9-1650996084562824.csv"
$ws2.Range("B5").Value = "OB-16509960851632621.csv"
$ws2.Range("B6").Value = "TB-16509960869792967.csv"
$ws2.Range("B7").Value = "TB-16509960865792928.csv"
$ws2.Range("B8").Value = "OB-16509960860112572.csv"
$ws2.Range("B9").Value = "ZB-match_3-16509960845868242.csv"
$ws2.Range("B10").Value = "TB-16509960866352944.csv"

# Sheet4: TOL_TO
$ws4 = $wb.Worksheets.Item("TOL_TO-16509960870432591")
$ws4.Range("B2").Value = "MM_stims-16509960870112906.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996086995259.csv"
$ws4.Range("B4").Value = "MM_stims-1650996087027291.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960870112906.csv"
$ws4.Range("B6").Value = "MM_stims-16509960870432591.csv"
$ws4.Range("B7").Value = "ZM_stims-1650996087027291.csv"

# Sheet5: vSAT_TO
$ws5 = $wb.Worksheets.Item("vSAT_TO-16509960871165123")
$ws5.Range("B2").Value = "SAT_stims-165099608705926.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509960870752938.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509960870919669.csv"
$ws5.Range("B5").Value = "SAT_stims-16509960870432591.csv"
